$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Winner") before the existing "Phase" column.
# This shifts the former B:E (Phase, First Gameweek, Last Gameweek,
# Number of Gameweeks) over to C:F.
$ws.Columns.Item(2).Insert()

# Header
$ws.Range("B1").Value = "Winner"

# MoTM winner per gameweek-phase row
$ws.Range("B2").Value  = "MB FC"
$ws.Range("B3").Value  = "ChicagoFire"
$ws.Range("B4").Value  = "Berba Juniors"
$ws.Range("B5").Value  = "Saka Souffle"
$ws.Range("B6").Value  = "Minnows FC"
$ws.Range("B7").Value  = "Gabi-Gabi-Gabagool"
$ws.Range("B8").Value  = "Momoney"
$ws.Range("B9").Value  = "TBD"
$ws.Range("B10").Value = "TBD"
$ws.Range("B11").Value = "TBD"
$ws.Range("B12").Value = "None"

# Match the new column's best-fit width (approx. 19.57 chars in the
# saved file).
$ws.Columns.Item(2).ColumnWidth = 18.71

# Restore the active selection used after entering this data.
$ws.Range("D17").Select()
